# Scheduled-runner update: refresh computed Leve profit figures
# (currentAveragePrice / profit columns H..N) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 474
$ws.Range("J55").Value = 550.17645
$ws.Range("L55").Value = 550.17645
$ws.Range("N55").Value = -978.17645

$ws.Range("H69").Value = 7267.8096
$ws.Range("I69").Value = 4970.3335
$ws.Range("K69").Value = 14911.0005
$ws.Range("M69").Value = -14037.0005

$ws.Range("H72").Value = 7267.8096
$ws.Range("I72").Value = 4970.3335
$ws.Range("K72").Value = 44733.0015
$ws.Range("M72").Value = -40365.0015

$ws.Range("H98").Value = 278356.22
$ws.Range("I98").Value = 1429.7333
$ws.Range("J98").Value = 871770.1
$ws.Range("K98").Value = 1429.7333
$ws.Range("L98").Value = 871770.1
$ws.Range("M98").Value = 68.2666999999999
$ws.Range("N98").Value = -874766.1

$ws.Range("H122").Value = 278356.22
$ws.Range("I122").Value = 1429.7333
$ws.Range("J122").Value = 871770.1
$ws.Range("K122").Value = 4289.199900000001
$ws.Range("L122").Value = 2615310.3
$ws.Range("M122").Value = -1839.199900000001
$ws.Range("N122").Value = -2620210.3

$ws.Range("H129").Value = 921
$ws.Range("J129").Value = 2825
$ws.Range("L129").Value = 8475
$ws.Range("N129").Value = -18475

$ws.Range("H132").Value = 1538.0435
$ws.Range("I132").Value = 1213.8948
$ws.Range("K132").Value = 3641.6844
$ws.Range("M132").Value = -1111.6844

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 2008999.8
$ws.Range("I8").Value = 5007499.5
$ws.Range("J8").Value = 10000
$ws.Range("K8").Value = 5007499.5
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = -5007355.5
$ws.Range("N8").Value = -10288

$ws.Range("H13").Value = 3339334.8
$ws.Range("J13").Value = 9002
$ws.Range("L13").Value = 9002
$ws.Range("N13").Value = -9290

$ws.Range("H56").Value = 9998.5
$ws.Range("I56").Value = 9998.5
$ws.Range("K56").Value = 9998.5
$ws.Range("M56").Value = -9256.5

$ws.Range("H61").Value = 5183.087
$ws.Range("I61").Value = 3509.95
$ws.Range("K61").Value = 3509.95
$ws.Range("M61").Value = -3297.95

$ws.Range("H98").Value = 30355
$ws.Range("J98").Value = 30355
$ws.Range("L98").Value = 30355
$ws.Range("N98").Value = -36345

$ws.Range("H136").Value = 5183.087
$ws.Range("I136").Value = 3509.95
$ws.Range("K136").Value = 10529.85
$ws.Range("M136").Value = -7979.849999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 56658
$ws.Range("J21").Value = 56658
$ws.Range("L21").Value = 56658
$ws.Range("N21").Value = -57130

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H105").Value = 14058.652
$ws.Range("I105").Value = 15737.929
$ws.Range("J105").Value = 11446.444
$ws.Range("K105").Value = 15737.929
$ws.Range("L105").Value = 11446.444
$ws.Range("M105").Value = -13990.929
$ws.Range("N105").Value = -14940.444

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 150.5238
$ws.Range("I7").Value = 69.75
$ws.Range("J7").Value = 258.22223
$ws.Range("K7").Value = 69.75
$ws.Range("L7").Value = 258.22223
$ws.Range("M7").Value = 43.25
$ws.Range("N7").Value = -484.22223

$ws.Range("H22").Value = 2798.889
$ws.Range("I22").Value = 1499.6666
$ws.Range("J22").Value = 3448.5
$ws.Range("K22").Value = 1499.6666
$ws.Range("L22").Value = 3448.5
$ws.Range("M22").Value = -1149.6666
$ws.Range("N22").Value = -4148.5

$ws.Range("H31").Value = 153561.14
$ws.Range("J31").Value = 153561.14
$ws.Range("L31").Value = 153561.14
$ws.Range("N31").Value = -154151.14

$ws.Range("H34").Value = 153561.14
$ws.Range("J34").Value = 153561.14
$ws.Range("L34").Value = 153561.14
$ws.Range("N34").Value = -153965.14

$ws.Range("H132").Value = 2945.027
$ws.Range("I132").Value = 2149.65
$ws.Range("J132").Value = 3880.7646
$ws.Range("K132").Value = 6448.950000000001
$ws.Range("L132").Value = 11642.2938
$ws.Range("M132").Value = -3918.950000000001
$ws.Range("N132").Value = -16702.2938

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 247.375
$ws.Range("J2").Value = 307.25
$ws.Range("L2").Value = 1843.5
$ws.Range("N2").Value = -2069.5

$ws.Range("H55").Value = 1376.9333
$ws.Range("I55").Value = 1229.5
$ws.Range("J55").Value = 1966.6666
$ws.Range("K55").Value = 3688.5
$ws.Range("L55").Value = 5899.9998
$ws.Range("M55").Value = -3511.5
$ws.Range("N55").Value = -6253.9998

$ws.Range("H138").Value = 1523.8
$ws.Range("I138").Value = 1523.8
$ws.Range("K138").Value = 4571.4
$ws.Range("M138").Value = 568.6000000000004

$ws.Range("H140").Value = 2455.4614
$ws.Range("I140").Value = 2110.5833
$ws.Range("K140").Value = 6331.749899999999
$ws.Range("M140").Value = -1151.749899999999

$ws.Range("H141").Value = 6076.9375
$ws.Range("J141").Value = 13626.6
$ws.Range("L141").Value = 40879.8
$ws.Range("N141").Value = -51239.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 37702.6
$ws.Range("J42").Value = 40083.668
$ws.Range("L42").Value = 40083.668
$ws.Range("N42").Value = -41053.668

$ws.Range("H113").Value = 2570.7576
$ws.Range("I113").Value = 2106
$ws.Range("J113").Value = 3500.2727
$ws.Range("K113").Value = 2106
$ws.Range("L113").Value = 3500.2727
$ws.Range("M113").Value = 64
$ws.Range("N113").Value = -7840.2727

$ws.Range("H115").Value = 37702.6
$ws.Range("J115").Value = 40083.668
$ws.Range("L115").Value = 40083.668
$ws.Range("N115").Value = -42433.668

$ws.Range("H126").Value = 4807
$ws.Range("I126").Value = 3473.6667
$ws.Range("J126").Value = 6140.3335
$ws.Range("K126").Value = 10421.0001
$ws.Range("L126").Value = 18421.0005
$ws.Range("M126").Value = -7951.000100000001
$ws.Range("N126").Value = -23361.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 999.5
$ws.Range("J12").Value = 999.5
$ws.Range("L12").Value = 999.5
$ws.Range("N12").Value = -1339.5

$ws.Range("H22").Value = 3240.4546
$ws.Range("I22").Value = 1949.5
$ws.Range("J22").Value = 4316.25
$ws.Range("K22").Value = 1949.5
$ws.Range("L22").Value = 4316.25
$ws.Range("M22").Value = -1654.5
$ws.Range("N22").Value = -4906.25

$ws.Range("H27").Value = 3240.4546
$ws.Range("I27").Value = 1949.5
$ws.Range("J27").Value = 4316.25
$ws.Range("K27").Value = 1949.5
$ws.Range("L27").Value = 4316.25
$ws.Range("M27").Value = -1842.5
$ws.Range("N27").Value = -4530.25

$ws.Range("H36").Value = 79000
$ws.Range("J36").Value = 79000
$ws.Range("L36").Value = 79000
$ws.Range("N36").Value = -80124

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1010000
$ws.Range("J8").Value = 1010000
$ws.Range("L8").Value = 1010000
$ws.Range("N8").Value = -1010280

$ws.Range("H13").Value = 799
$ws.Range("J13").Value = 799
$ws.Range("L13").Value = 799
$ws.Range("N13").Value = -1079

$ws.Range("H81").Value = 5226.8125
$ws.Range("J81").Value = 6074.75
$ws.Range("L81").Value = 12149.5
$ws.Range("N81").Value = -14271.5

$ws.Range("H84").Value = 5226.8125
$ws.Range("J84").Value = 6074.75
$ws.Range("L84").Value = 60747.5
$ws.Range("N84").Value = -71355.5

$ws.Range("H97").Value = 2500
$ws.Range("J97").Value = 2500
$ws.Range("L97").Value = 2500
$ws.Range("N97").Value = -4482

$ws.Range("H107").Value = 1497.091
$ws.Range("I107").Value = 1709.3334
$ws.Range("J107").Value = 542
$ws.Range("K107").Value = 5128.0002
$ws.Range("L107").Value = 1626
$ws.Range("M107").Value = -3208.0002
$ws.Range("N107").Value = -5466
